$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain decimal numbers need to be forced to
# text format first, otherwise Excel auto-converts them to numbers and silently
# drops significant trailing zeros (e.g. "300.84" vs "0.520" -> 0.52).
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "300.84"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "98.17"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.520"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.517"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.88"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0792"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.790"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.24"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.13"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "68.40"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "238.57"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.21"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "25.19"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "167.46"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.03"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "32.97"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.00"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.15"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "18.23"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0689"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.20"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "17.56"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.84"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "54.52"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.54"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.80"

$ws.Range("D2").Value = "43.161.52"
$ws.Range("E2").Value = "  -0.33%  "
$ws.Range("D3").Value = "2.310.04"
$ws.Range("E3").Value = "  +0.04%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("E5").Value = "  -0.46%  "
$ws.Range("E6").Value = "  -2.70%  "
$ws.Range("E7").Value = "  +3.37%  "
$ws.Range("E8").Value = "  -0.09%  "
$ws.Range("E9").Value = "  +0.07%  "
$ws.Range("E10").Value = "  -2.22%  "
$ws.Range("E11").Value = "  -0.36%  "
$ws.Range("E12").Value = "  +0.22%  "
$ws.Range("E13").Value = "  -2.88%  "
$ws.Range("E14").Value = "  -1.27%  "
$ws.Range("D15").Value = "2.669.56"
$ws.Range("E15").Value = "  -0.08%  "
$ws.Range("D16").Value = "2.261.13"
$ws.Range("E16").Value = "  -3.08%  "
$ws.Range("E17").Value = "  -1.87%  "
$ws.Range("D18").Value = "43.076.55"
$ws.Range("E18").Value = "  -0.21%  "
$ws.Range("E19").Value = "  +4.09%  "
$ws.Range("D20").Value = "0.0₃0911"
$ws.Range("E20").Value = "  +0.06%  "
$ws.Range("E21").Value = "  -1.68%  "
$ws.Range("E22").Value = "  +0.29%  "
$ws.Range("E23").Value = "  +0.65%  "
$ws.Range("E24").Value = "  -1.84%  "
$ws.Range("E25").Value = "  -0.07%  "
$ws.Range("E26").Value = "  -1.93%  "
$ws.Range("E27").Value = "  +0.04%  "
$ws.Range("E28").Value = "  -0.59%  "
$ws.Range("E29").Value = "  -0.44%  "
$ws.Range("E30").Value = "  -2.14%  "
$ws.Range("E31").Value = "  -5.90%  "
$ws.Range("B32").Value = "FirstDigitalUSD"
$ws.Range("C32").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("E32").Value = "  -0.05%  "
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("E33").Value = "  +2.27%  "
$ws.Range("E34").Value = "  +2.16%  "
$ws.Range("E35").Value = "  +1.69%  "
$ws.Range("E36").Value = "  -0.38%  "
$ws.Range("E37").Value = "  -1.32%  "
$ws.Range("E38").Value = "  +1.33%  "
$ws.Range("E39").Value = "  -0.28%  "
$ws.Range("E40").Value = "  +1.27%  "
$ws.Range("E41").Value = "  -3.16%  "
$ws.Range("D42").Value = "2.017.59"
$ws.Range("E42").Value = "  +1.30%  "
$ws.Range("E43").Value = "  -1.07%  "
$ws.Range("E44").Value = "  -6.27%  "
$ws.Range("E45").Value = "  +0.26%  "
$ws.Range("E46").Value = "  -1.13%  "
$ws.Range("E47").Value = "  -2.56%  "
$ws.Range("E48").Value = "  -2.46%  "
$ws.Range("D49").Value = "2.534.49"
$ws.Range("E49").Value = "  -0.20%  "
$ws.Range("E50").Value = "  -1.52%  "
$ws.Range("B51").Value = "HuobiToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("E51").Value = "  +14.65%  "
